$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1831
$ws1.Range("F6").Value = 89
$ws1.Range("F10").Value = 75
$ws1.Range("F12").Value = 5257
$ws1.Range("F14").Value = 870
$ws1.Range("F16").Value = 2336
$ws1.Range("F19").Value = 2178

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1831
$ws4.Range("F6").Value = 89
$ws4.Range("F10").Value = 75
$ws4.Range("F12").Value = 5257
$ws4.Range("F16").Value = 870
$ws4.Range("F18").Value = 2336
$ws4.Range("F22").Value = 2178
